# Add "Maven Lifecycle default bindings" entries to Sheet1 (git + maven lib commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16 - Clean lifecycle bindings
$ws.Cells.Item(16, 1).Value = "Maven"
$ws.Cells.Item(16, 2).Value = "LifeCycle default binding to plugins (1)"
$ws.Cells.Item(16, 3).Value = @"
Clean Lifecycle Bindings (Plugin)
* clean//clean:clean
"@

# Row 17 - Default lifecycle bindings for ejb/ejb3/jar/par/rar/war packaging
$ws.Cells.Item(17, 1).Value = "Maven"
$ws.Cells.Item(17, 2).Value = "LifeCycle default binding to plugins (2)"
$ws.Cells.Item(17, 3).Value = @"
Default Lifecycle Bindings - Packaging ejb / ejb3 / jar / par / rar / war (Plugin)
* process-resources//resources:resources
* compile//compiler:compile
* process-test-resources//resources:testResources
* test-compile//compiler:testCompile
* test//surefire:test
* package//ejb:ejb or ejb3:ejb3 or jar:jar or par:par or rar:rar or war:war
* install//install:install
* deploy//deploy:deploy
"@

# Row 18 - Default lifecycle bindings for ear packaging
$ws.Cells.Item(18, 1).Value = "Maven"
$ws.Cells.Item(18, 2).Value = "LifeCycle default binding to plugins (3)"
$ws.Cells.Item(18, 3).Value = @"
Default Lifecycle Bindings - Packaging ear (Plugin)
* generate-resources//ear:generate-application-xml
* process-resources//resources:resources
* package//ear:ear
* install//install:install
* deploy//deploy:deploy
"@

# Row 19 - Default lifecycle bindings for maven-plugin packaging
$ws.Cells.Item(19, 1).Value = "Maven"
$ws.Cells.Item(19, 2).Value = "LifeCycle default binding to plugins (4)"
$ws.Cells.Item(19, 3).Value = @"
Default Lifecycle Bindings - Packaging maven-plugin (Plugin)
* generate-resources//plugin:descriptor
* process-resources//resources:resources
* compile//compiler:compile
* process-test-resources//resources:testResources
* test-compile//compiler:testCompile
* test//surefire:test
* package//jar:jar and plugin:addPluginArtifactMetadata
* install//install:install
* deploy//deploy:deploy
"@

# Row 20 - Default lifecycle bindings for pom packaging
$ws.Cells.Item(20, 1).Value = "Maven"
$ws.Cells.Item(20, 2).Value = "LifeCycle default binding to plugins (5)"
$ws.Cells.Item(20, 3).Value = @"
Default Lifecycle Bindings - Packaging pom (Plugin)
* package//site:attach-descriptor
* install//install:install
* deploy//deploy:deploy
"@

# Row 21 - Site lifecycle bindings
$ws.Cells.Item(21, 1).Value = "Maven"
$ws.Cells.Item(21, 2).Value = "LifeCycle default binding to plugins (6)"
$ws.Cells.Item(21, 3).Value = @"
Site Lifecycle Bindings (Plugin)
* site//site:site
* site-deploy//site:deploy
"@

# Column B has no column-level style (unlike A & C), so new cells in B17:B21
# don't pick up the wrap-text style automatically the way B16 (pre-existing
# blank cell) does. Clone B16's format onto them.
$ws.Cells.Item(16, 2).Copy()
$ws.Range("B17:B21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the sheet's standard row height to the new rows (keeps them in
# line with the existing 36pt data rows rather than an autofit value).
$ws.Range("A16:C21").RowHeight = 36

# Match the author's final selection/scroll state shown in the diff.
$ws.Range("A17:A21").Select()
